$wb = $excel.ActiveWorkbook

# Both "Score" and "Rank" sheets reference the same two shared strings in
# B2/C2 (the ratio and its description). The shared string table text is
# replaced in place, so update B2/C2 on every worksheet that uses them so
# they end up pointing at the same (new) shared text.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Range("B2").Value = "0.21:0.07:0.25:0.3:0.17"
    $sheet.Range("C2").Value = "0.21T, 0.07RR, 0.25Env, 0.3Econ, 0.17S"
}

# Only the "Score" sheet's results (D2:F2) actually change value.
$ws = $wb.Worksheets.Item("Score")
$ws.Range("D2").Value = 0.1755158851424245
$ws.Range("E2").Value = 0.8578473045443887
$ws.Range("F2").Value = 0.05851766611794292
